$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reorderings (column A) ---
$ws.Range("A144").Value = "Zambia"
$ws.Range("A145").Value = "Puerto Rico"
$ws.Range("A147").Value = "Islas Caimanes"
$ws.Range("A148").Value = "Bermudas"
$ws.Range("A149").Value = "Guam"
$ws.Range("A150").Value = "Eritrea"
$ws.Range("A162").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("A163").Value = "Libia"
$ws.Range("A164").Value = "Nueva Caledonia"
$ws.Range("A181").Value = "Sudan"
$ws.Range("A182").Value = "Angola"
$ws.Range("A183").Value = "Liberia"
$ws.Range("A184").Value = "Suazilandia"
$ws.Range("A185").Value = "Republica del Chad"
$ws.Range("A187").Value = "Nepal"
$ws.Range("A188").Value = "Zimbabue"
$ws.Range("A192").Value = "Somalia"
$ws.Range("A193").Value = "San Vicente y las Granadinas"
$ws.Range("A194").Value = "Cabo Verde"
$ws.Range("A200").Value = "Malaui"
$ws.Range("A202").Value = "Belice"
$ws.Range("A205").Value = "Islas Virgenes Britanicas"
$ws.Range("A206").Value = "Anguila"
$ws.Range("A207").Value = "Burundi"
$ws.Range("A210").Value = "Papua Nueva Guinea"
$ws.Range("A211").Value = "Timor Oriental"

# --- Numeric data updates (columns B-H) ---
$ws.Range("B4").Value = 306854
$ws.Range("C4").Value = 29693
$ws.Range("E4").Value = 283818
$ws.Range("G4").Value = 946
$ws.Range("H4").Value = 8350
$ws.Range("B5").Value = 124870
$ws.Range("C5").Value = 5671
$ws.Range("E5").Value = 78833
$ws.Range("G5").Value = 620
$ws.Range("H5").Value = 11818
$ws.Range("B7").Value = 95614
$ws.Range("C7").Value = 4455
$ws.Range("E7").Value = 67787
$ws.Range("G7").Value = 152
$ws.Range("H7").Value = 1427
$ws.Range("B16").Value = 13912
$ws.Range("C16").Value = 1537
$ws.Range("E16").Value = 11086
$ws.Range("E90").Value = 215
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 19
$ws.Range("B102").Value = 217
$ws.Range("C102").Value = 23
$ws.Range("E102").Value = 195
$ws.Range("B104").Value = 209
$ws.Range("E104").Value = 180
$ws.Range("D144").Value = 2
$ws.Range("H144").Value = 1
$ws.Range("D145").Value = 1
$ws.Range("H145").Value = 2
$ws.Range("C147").Value = 6
$ws.Range("D147").Value = 1
$ws.Range("E147").Value = 33
$ws.Range("H147").Value = 1
$ws.Range("B148").Value = 35
$ws.Range("D148").Value = 14
$ws.Range("E148").Value = 21
$ws.Range("H148").Value = 0
$ws.Range("B149").Value = 32
$ws.Range("C149").Value = 0
$ws.Range("E149").Value = 31
$ws.Range("H149").Value = 1
$ws.Range("C150").Value = 7
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 29
$ws.Range("H150").Value = 0
$ws.Range("B162").Value = 17
$ws.Range("D162").Value = 0
$ws.Range("E163").Value = 16
$ws.Range("H163").Value = 1
$ws.Range("D164").Value = 1
$ws.Range("H164").Value = 0
$ws.Range("C181").Value = 0
$ws.Range("C182").Value = 2
$ws.Range("D182").Value = 2
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 2
$ws.Range("C183").Value = 3
$ws.Range("D183").Value = 3
$ws.Range("G183").Value = 1
$ws.Range("H183").Value = 1
$ws.Range("C184").Value = 0
$ws.Range("C185").Value = 1
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 1
$ws.Range("H187").Value = 0
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 0
$ws.Range("H188").Value = 1
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 1
$ws.Range("H193").Value = 0
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 0
$ws.Range("H194").Value = 1
$ws.Range("C200").Value = 1
$ws.Range("C202").Value = 0

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 23:22"
